$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated dSF (column F) values re-pulled from the source data.
$updates = @{
    2  = -1
    4  = 3
    6  = 1
    7  = -1
    9  = 0
    11 = -3
    16 = -5
    20 = -1
    21 = -7
    23 = -1
    27 = -9
    28 = -4
    31 = -2
    32 = -8
    33 = -1
    34 = 0
    35 = 2
    37 = 0
    46 = -2
    47 = -4
    48 = -7
    49 = -1
    57 = -1
    58 = -2
    60 = -2
    65 = 1
    66 = -1
    67 = -1
    68 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
